# Update timing_results sheet:
# - Insert a new "pyboolnet" column before the current Pint column,
#   shifting CABEAN/PyBoolNet_Asp/stable_motifs_new/bioLQM/boolsim/sm_jgtz/PyBoolNet/model
#   columns over, renaming several of them, and dropping one column (J) overall.
# - Replace row 2 timing values with the latest benchmark run results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old J column entirely (will no longer be used; new layout is A:I only)
$ws.Range("J1:J2").Clear()

# --- Header row (row 1) ---
$ws.Range("B1").Value = "pyboolnet"
$ws.Range("C1").Value = "Pint"
$ws.Range("D1").Value = "boolsim"
$ws.Range("E1").Value = "stablemotifs2013"
$ws.Range("F1").Value = "bioLQM"
$ws.Range("G1").Value = "CABEAN"
$ws.Range("H1").Value = "pystablemotifs"
$ws.Range("I1").Value = "model"

# Copy header style (bold/centered/bordered) from B1 onto the new I1 cell,
# matching the style already used across the header row.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data row (row 2) ---
$ws.Range("B2").Value = 0.5727500239991059
$ws.Range("C2").Value = 0.3547523589986668
$ws.Range("D2").Value = 0.2571681939989503
$ws.Range("E2").Value = 59.29691573200034
$ws.Range("F2").Value = 0.1880386279990489
$ws.Range("G2").Value = 0.1789524010018795
$ws.Range("H2").Value = 0.9352700599993113
$ws.Range("I2").Value = "phase_switch"

$ws.UsedRange | Out-Null
